$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.365400195121765
$ws.Range("B1").Value = 1.768108248710632
$ws.Range("C1").Value = 1.531245470046997
$ws.Range("D1").Value = 2.262045621871948
$ws.Range("E1").Value = 3.58557915687561
